# "data helper clas has beeen modified"
# The dataForLogin sheet gains a third column ("login") that records
# whether each UserName/Password pair resulted in a successful login,
# the "hareg" row is dropped, and a couple of rows are adjusted/added
# to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataForLogin")

# New header for column C
$ws.Range("C1").Value = "login"

# Rows 2-9: keep existing UserName/Password, add the login result column
$ws.Range("C2").Value = "success"
$ws.Range("C3").Value = "success"
$ws.Range("C4").Value = "fail"
$ws.Range("C5").Value = "fail"
$ws.Range("C6").Value = "fail"
$ws.Range("C7").Value = "fail"
$ws.Range("C8").Value = "fail"
$ws.Range("C9").Value = "fail"

# Row 10 ("hareg","hareg") is removed - clear it out
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "fail"

# Row 11 keeps "salima" in A, but B is cleared
$ws.Range("A11").Value = "salima"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "fail"

# New row 12 carries a stray "a" in column B
$ws.Range("A12").Value = ""
$ws.Range("B12").Value = "a"
$ws.Range("C12").Value = "fail"

# Update selection to match the author's final cursor position
$ws.Range("B10").Select()
